# "add a false case" -- append a new row 8 to Sheet1: A8="abd", B8="2018",
# C8="False", all three cells sharing the text-formatted style already used
# by (empty) C8, plus move the active selection to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C8 already carries the row's "text" cell style (General/Text number format
# applied per-cell). Clone that style onto A8:B8 first so every cell in the
# new row matches formatting-wise, without minting extra style records.
$ws.Range("C8").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A8").Value = "abd"
$ws.Range("B8").Value = "2018"

# Plain `Value = "False"` gets auto-coerced to the Boolean FALSE by the
# engine's cell-input parser (regardless of the cell's Text format), which
# is not what we want here -- the source data is the literal word "False".
# Route the literal through a formula result first (a formula's text result
# is never re-interpreted as a boolean) and paste only the resulting value
# onto C8, leaving C8's existing style untouched.
$ws.Range("Z1").Formula = '="False"'
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial(-4163)      # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("F11").Select()
